# Append new otolith records (rows 175-215) to Sheet1, mirroring the rows
# already present in the sheet (Species column in italics via style "s=2",
# plain-number Trawl/No.Extracted columns, text Well No./Plate No. columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Trawl, Species, Length, No. Extracted, Well No., Plate No.
$rows = @(
    @(6, "Macroparalepis affinis", "96mm SL", 2, "D1", "TC3"),
    @(6, "Macroparalepis affinis", "127mm SL", 1, "D2", "TC3"),
    @(6, "Macroparalepis affinis", "96mm SL", 2, "D3", "TC3"),
    @(6, "Macroparalepis affinis", "NR", 2, "D4", "TC3"),
    @(6, "Macroparalepis affinis", "92mm SL", 2, "D5", "TC3"),
    @(6, "Macroparalepis affinis", "98mm SL", 2, "D6", "TC3"),
    @(6, "Macroparalepis affinis", "92mm SL", 2, "D7", "TC3"),
    @(6, "Macroparalepis affinis", "93mm SL", 2, "D8", "TC3"),
    @(6, "Macroparalepis affinis", "112mm SL", 2, "D9", "TC3"),
    @(6, "Macroparalepis affinis", "94mm SL", 2, "D10", "TC3"),
    @(6, "Macroparalepis affinis", "113mm SL", 2, "D11", "TC3"),
    @(6, "Macroparalepis affinis", "82mm SL", 2, "D12", "TC3"),
    @(6, "Macroparalepis affinis", "122mm SL", 2, "F1", "TC3"),
    @(6, "Macroparalepis affinis", "111mm SL", 2, "F2", "TC3"),
    @(6, "Macroparalepis affinis", "91mm SL", 2, "F3", "TC3"),
    @(6, "Macroparalepis affinis", "103mm SL", 2, "F4", "TC3"),
    @(6, "Macroparalepis affinis", "95mm SL", 2, "F5", "TC3"),
    @(6, "Macroparalepis affinis", "118mm SL", 2, "F6", "TC3"),
    @(6, "Macroparalepis affinis", "115mm SL", 2, "F7", "TC3"),
    @(6, "Macroparalepis affinis", "136mm SL", 2, "F8", "TC3"),
    @(8, "Argyropelecus olfersii", "NR", 2, "E4", "TC3"),
    @(8, "Trachurus trachurus", "63mm SL", 2, "B1", "TC3"),
    @(8, "Trachurus trachurus", "64mm SL", 2, "B2", "TC3"),
    @(4, "Maurolicus muelleri", "NR", 2, "F1", "TC2"),
    @(4, "Maurolicus muelleri", "42mm SL", 2, "F2", "TC2"),
    @(4, "Maurolicus muelleri", "NR", 2, "F3", "TC2"),
    @(4, "Maurolicus muelleri", "45mm SL", 2, "F4", "TC2"),
    @(4, "Maurolicus muelleri", "47mm SL", 2, "F5", "TC2"),
    @(4, "Maurolicus muelleri", "46mm SL", 2, "F6", "TC2"),
    @(4, "Maurolicus muelleri", "NR", 2, "F7", "TC2"),
    @(4, "Maurolicus muelleri", "47mm SL", 2, "F8", "TC2"),
    @(4, "Maurolicus muelleri", "45mm SL", 2, "F9", "TC2"),
    @(4, "Maurolicus muelleri", "49mm SL", 2, "F10", "TC2"),
    @(4, "Maurolicus muelleri", "47mm SL", 2, "F11", "TC2"),
    @(4, "Maurolicus muelleri", "45mm SL", 2, "F12", "TC2"),
    @(4, "Maurolicus muelleri", "47mm SL", 2, "G1", "TC2"),
    @(4, "Maurolicus muelleri", "49mm SL", 2, "G2", "TC2"),
    @(4, "Maurolicus muelleri", "NR", 2, "G3", "TC2"),
    @(4, "Maurolicus muelleri", "46mm SL", 2, "G4", "TC2"),
    @(4, "Maurolicus muelleri", "44mm SL", 2, "G5", "TC2"),
    @(4, "Maurolicus muelleri", "49mm SL", 2, "G6", "TC2")
)

$startRow = 175
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i

    # The original workbook's shared-string table grew with "44mm SL" (row 214's
    # length) entering before "49mm SL" (first seen on row 207), even though row
    # 207 precedes row 214 in the sheet. Reproduce that exact insertion order by
    # poking row 214's Length cell a touch early, right before row 207 is written.
    if ($r -eq 207) {
        $ws.Cells.Item(214, 3).Value = $rows[214 - $startRow][2]
    }

    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 2).Font.Italic = $true

    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

$lastRow = $startRow + $rows.Count - 1
$ws.Range("A$lastRow").Offset(1, 0).Select()
